$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Add the new header / label cells for the "Room" block (rows 14-16) ---
$ws.Range("B14").Value = "Room"
$ws.Range("H14").Value = "Roomuse"

$ws.Range("B15").Value = "คอลัมน์1"
$ws.Range("C15").Value = "คอลัมน์2"
$ws.Range("H15").Value = "คอลัมน์1"
$ws.Range("I15").Value = "คอลัมน์2"
$ws.Range("J15").Value = "คอลัมน์3"
$ws.Range("K15").Value = "คอลัมน์4"

$ws.Range("B16").Value = "RoomNum"
$ws.Range("C16").Value = "RoomType"
$ws.Range("H16").Value = "RoomNum"
$ws.Range("I16").Value = "Per_ID"
$ws.Range("J16").Value = "UseTime"
$ws.Range("K16").Value = "LeaveTime"

# --- 2. Rebuild the existing tables so the engine can append new ones -------
# The workbook's table id sequence has a gap (1,2,4,5,6,7 -- id 3 is free
# because a table was removed earlier in the workbook's history). The
# ListObjects collection can only append a new table while the id sequence
# it already has is contiguous, so temporarily "unlist" (which keeps the
# cell data, it just drops the table/ListObject wrapper) the four tables
# that sit after the gap and immediately recreate them identically. This
# leaves the sheet content totally unchanged while clearing the gap.
$ws.ListObjects.Item(6).Unlist()   # Table7  H11:J12
$ws.ListObjects.Item(5).Unlist()   # Table6  B11:D12
$ws.ListObjects.Item(4).Unlist()   # Table5  H7:J8
$ws.ListObjects.Item(3).Unlist()   # Table4  B7:F8

# --- 3. Re-create the four original tables plus the two brand-new ones -----
# in the exact order that reproduces the target id/name scheme:
#   Table1 (untouched), Table2 (untouched), Table3 (new, B15:C16),
#   Table4, Table5, Table6, Table7 (restored), Table8 (new, H15:K16)
$ws.ListObjects.Add(1, $ws.Range("B15:C16"), $null, 1) | Out-Null
$ws.ListObjects.Add(1, $ws.Range("B7:F8"), $null, 1) | Out-Null
$ws.ListObjects.Add(1, $ws.Range("H7:J8"), $null, 1) | Out-Null
$ws.ListObjects.Add(1, $ws.Range("B11:D12"), $null, 1) | Out-Null
$ws.ListObjects.Add(1, $ws.Range("H11:J12"), $null, 1) | Out-Null
$ws.ListObjects.Add(1, $ws.Range("H15:K16"), $null, 1) | Out-Null

# --- 4. Match the selection left behind by the author's edit ---------------
$ws.Range("H15:K16").Select()
